# Update imputed values in columns D and E of Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 15.06879999999999
$ws.Range("D8").Value = -8.581700000000001
$ws.Range("D10").Value = -9.459599999999991
$ws.Range("D12").Value = -7.121699999999997
$ws.Range("E15").Value = 16.0118
$ws.Range("D18").Value = -8.523199999999999
$ws.Range("E18").Value = 16.48570000000001
$ws.Range("E20").Value = 15.9051
$ws.Range("E29").Value = 17.06870000000001
$ws.Range("E30").Value = 15.51389999999999
$ws.Range("E31").Value = 15.98309999999999
$ws.Range("D37").Value = -7.750399999999994
$ws.Range("E40").Value = 17.05060000000002
$ws.Range("E50").Value = 16.3091
$ws.Range("D55").Value = -8.401999999999999
$ws.Range("D68").Value = -7.025399999999998
$ws.Range("E68").Value = 18.1281
$ws.Range("E76").Value = 16.25209999999998
$ws.Range("D77").Value = -5.748500000000002
$ws.Range("D78").Value = -7.518700000000004
$ws.Range("D81").Value = -7.641499999999997
$ws.Range("D82").Value = -8.214600000000001
$ws.Range("E87").Value = 16.2987
$ws.Range("E88").Value = 16.3699
$ws.Range("E96").Value = 16.11639999999999
$ws.Range("E98").Value = 15.3486
$ws.Range("E101").Value = 16.71900000000002
$ws.Range("E102").Value = 16.6573
